$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4784.357
$ws.Range("J17").Value = 5298.48
$ws.Range("L17").Value = 15895.44
$ws.Range("N17").Value = -16231.44
$ws.Range("H28").Value = 2474.5264
$ws.Range("I28").Value = 501.35715
$ws.Range("K28").Value = 501.35715
$ws.Range("M28").Value = -16.35714999999999
$ws.Range("H70").Value = 5600.625
$ws.Range("I70").Value = 1376.8572
$ws.Range("J70").Value = 6321.756
$ws.Range("K70").Value = 4130.571599999999
$ws.Range("L70").Value = 18965.268
$ws.Range("M70").Value = -3860.571599999999
$ws.Range("N70").Value = -19505.268
$ws.Range("H73").Value = 5600.625
$ws.Range("I73").Value = 1376.8572
$ws.Range("J73").Value = 6321.756
$ws.Range("K73").Value = 4130.571599999999
$ws.Range("L73").Value = 18965.268
$ws.Range("M73").Value = -3194.571599999999
$ws.Range("N73").Value = -20837.268
$ws.Range("H127").Value = 1265.7273
$ws.Range("I127").Value = 892.3
$ws.Range("K127").Value = 2676.9
$ws.Range("M127").Value = 2283.1
$ws.Range("H132").Value = 32261804
$ws.Range("I132").Value = 40004172
$ws.Range("K132").Value = 120012516
$ws.Range("M132").Value = -120009986
$ws.Range("H135").Value = 890.7273
$ws.Range("I135").Value = 861.2353000000001
$ws.Range("K135").Value = 7751.117700000001
$ws.Range("M135").Value = -5216.117700000001
$ws.Range("H138").Value = 4374.9653
$ws.Range("I138").Value = 2916.3333
$ws.Range("J138").Value = 4755.478
$ws.Range("K138").Value = 8748.999899999999
$ws.Range("L138").Value = 14266.434
$ws.Range("M138").Value = -3608.999899999999
$ws.Range("N138").Value = -24546.434
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3969656.8
$ws.Range("I2").Value = 5051581.5
$ws.Range("K2").Value = 5051581.5
$ws.Range("M2").Value = -5051468.5
$ws.Range("H61").Value = 14648.889
$ws.Range("I61").Value = 14648.889
$ws.Range("K61").Value = 14648.889
$ws.Range("M61").Value = -14436.889
$ws.Range("H97").Value = 1704172
$ws.Range("I97").Value = 2023610.8
$ws.Range("K97").Value = 2023610.8
$ws.Range("M97").Value = -2023114.8
$ws.Range("H102").Value = 3626638.8
$ws.Range("I102").Value = 6947340.5
$ws.Range("J102").Value = 4055.2727
$ws.Range("K102").Value = 6947340.5
$ws.Range("L102").Value = 4055.2727
$ws.Range("M102").Value = -6945718.5
$ws.Range("N102").Value = -7299.2727
$ws.Range("H116").Value = 3969656.8
$ws.Range("I116").Value = 5051581.5
$ws.Range("K116").Value = 5051581.5
$ws.Range("M116").Value = -5049287.5
$ws.Range("H136").Value = 14648.889
$ws.Range("I136").Value = 14648.889
$ws.Range("K136").Value = 43946.667
$ws.Range("M136").Value = -41396.667
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3969656.8
$ws.Range("I3").Value = 5051581.5
$ws.Range("K3").Value = 5051581.5
$ws.Range("M3").Value = -5051467.5
$ws.Range("H86").Value = 5891760.5
$ws.Range("I86").Value = 14307057
$ws.Range("K86").Value = 14307057
$ws.Range("M86").Value = -14305934
$ws.Range("H89").Value = 5891760.5
$ws.Range("I89").Value = 14307057
$ws.Range("K89").Value = 71535285
$ws.Range("M89").Value = -71529669
$ws.Range("H99").Value = 3325002.5
$ws.Range("I99").Value = 4466770
$ws.Range("K99").Value = 4466770
$ws.Range("M99").Value = -4465272
$ws.Range("H105").Value = 2500836.5
$ws.Range("I105").Value = 2841569.2
$ws.Range("K105").Value = 2841569.2
$ws.Range("M105").Value = -2839822.2
$ws.Range("H107").Value = 2859126.8
$ws.Range("I107").Value = 3968961.5
$ws.Range("K107").Value = 3968961.5
$ws.Range("M107").Value = -3967041.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 8
$ws.Range("I17").Value = 8
$ws.Range("K17").Value = 8
$ws.Range("M17").Value = 166
$ws.Range("H22").Value = 553
$ws.Range("I22").Value = 523.06665
$ws.Range("K22").Value = 523.06665
$ws.Range("M22").Value = -173.06665
$ws.Range("H58").Value = 2789.6667
$ws.Range("I58").Value = 2053.4167
$ws.Range("K58").Value = 2053.4167
$ws.Range("M58").Value = -1850.4167
$ws.Range("H59").Value = 26833.334
$ws.Range("J59").Value = 50500
$ws.Range("L59").Value = 50500
$ws.Range("N59").Value = -52790
$ws.Range("H62").Value = 3865.3333
$ws.Range("I62").Value = 3112.7144
$ws.Range("K62").Value = 3112.7144
$ws.Range("M62").Value = -2488.7144
$ws.Range("H65").Value = 3865.3333
$ws.Range("I65").Value = 3112.7144
$ws.Range("K65").Value = 15563.572
$ws.Range("M65").Value = -12443.572
$ws.Range("H70").Value = 10000
$ws.Range("J70").Value = 10000
$ws.Range("L70").Value = 10000
$ws.Range("N70").Value = -10630
$ws.Range("H73").Value = 10000
$ws.Range("J73").Value = 10000
$ws.Range("L73").Value = 10000
$ws.Range("N73").Value = -12184
$ws.Range("H86").Value = 12113.407
$ws.Range("I86").Value = 10267.357
$ws.Range("K86").Value = 10267.357
$ws.Range("M86").Value = -9144.357
$ws.Range("H89").Value = 12113.407
$ws.Range("I89").Value = 10267.357
$ws.Range("K89").Value = 51336.785
$ws.Range("M89").Value = -45720.785
$ws.Range("H93").Value = 31618.375
$ws.Range("J93").Value = 92399.5
$ws.Range("L93").Value = 92399.5
$ws.Range("N93").Value = -96143.5
$ws.Range("H105").Value = 2868.1304
$ws.Range("I105").Value = 2645.7058
$ws.Range("J105").Value = 3498.3333
$ws.Range("K105").Value = 2645.7058
$ws.Range("L105").Value = 3498.3333
$ws.Range("M105").Value = -898.7058000000002
$ws.Range("N105").Value = -6992.3333
$ws.Range("H122").Value = 2614.3157
$ws.Range("I122").Value = 2480.4
$ws.Range("K122").Value = 7441.200000000001
$ws.Range("M122").Value = -4991.200000000001
$ws.Range("H134").Value = 19214.154
$ws.Range("I134").Value = 35963.332
$ws.Range("J134").Value = 14189.4
$ws.Range("K134").Value = 107889.996
$ws.Range("L134").Value = 42568.2
$ws.Range("M134").Value = -105354.996
$ws.Range("N134").Value = -47638.2
$ws.Range("H136").Value = 2789.6667
$ws.Range("I136").Value = 2053.4167
$ws.Range("K136").Value = 6160.250100000001
$ws.Range("M136").Value = -3610.250100000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 933.55554
$ws.Range("I3").Value = 933.55554
$ws.Range("K3").Value = 2800.66662
$ws.Range("M3").Value = -2688.66662
$ws.Range("H11").Value = 7033.4165
$ws.Range("I11").Value = 2377.889
$ws.Range("K11").Value = 7133.667
$ws.Range("M11").Value = -6993.667
$ws.Range("H130").Value = 2199
$ws.Range("I130").Value = 997.5
$ws.Range("J130").Value = 3000
$ws.Range("K130").Value = 2992.5
$ws.Range("L130").Value = 9000
$ws.Range("M130").Value = 2027.5
$ws.Range("N130").Value = -19040
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4000229
$ws.Range("J102").Value = 1912698.4
$ws.Range("L102").Value = 1912698.4
$ws.Range("N102").Value = -1915942.4
$ws.Range("H113").Value = 10419134
$ws.Range("I113").Value = 20834228
$ws.Range("J113").Value = 4038.875
$ws.Range("K113").Value = 20834228
$ws.Range("L113").Value = 4038.875
$ws.Range("M113").Value = -20832058
$ws.Range("N113").Value = -8378.875
$ws.Range("H122").Value = 813870.4
$ws.Range("I122").Value = 894857.4
$ws.Range("K122").Value = 2684572.2
$ws.Range("M122").Value = -2682122.2
$ws.Range("H132").Value = 11139.556
$ws.Range("I132").Value = 9036.941000000001
$ws.Range("J132").Value = 14714
$ws.Range("K132").Value = 27110.823
$ws.Range("L132").Value = 44142
$ws.Range("M132").Value = -24580.823
$ws.Range("N132").Value = -49202
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 19000
$ws.Range("J108").Value = 19000
$ws.Range("L108").Value = 19000
$ws.Range("N108").Value = -26680
$ws.Range("H132").Value = 12616.685
$ws.Range("I132").Value = 14294.875
$ws.Range("J132").Value = 3666.3333
$ws.Range("K132").Value = 42884.625
$ws.Range("L132").Value = 10998.9999
$ws.Range("M132").Value = -40354.625
$ws.Range("N132").Value = -16058.9999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 166666670
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("H84").Value = 166666670
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("H96").Value = 1939.4
$ws.Range("H132").Value = 20638052
$ws.Range("I132").Value = 25006628
$ws.Range("J132").Value = 1222151
$ws.Range("K132").Value = 25006628
$ws.Range("L132").Value = 3666453
$ws.Range("M132").Value = -75017354
$ws.Range("N132").Value = -3671513

# Cell deletions (diff removes these cells entirely)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N81").ClearContents()
$ws.Range("N84").ClearContents()
